$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for all data rows (2-115) from 45178 to 45179
for ($r = 2; $r -le 115; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}

# Row 4 specific updates: Signalarter count (I4) and Alla arter count (Q4)
$ws.Range("I4").Value = 10
$ws.Range("Q4").Value = 14

# Row 4 Artnamn list (R4): insert "Rödgul trumpetsvamp" alphabetically
# between "Klippfrullania" and "Smal svampklubba"
$crlf = "`r`n"
$artnamn4 = @(
    "Entita",
    "Mindre hackspett",
    "Spillkråka",
    "Talltita",
    "Blåmossa",
    "Bronshjon",
    "Fällmossa",
    "Guldlockmossa",
    "Klippfrullania",
    "Rödgul trumpetsvamp",
    "Smal svampklubba",
    "Stor revmossa",
    "Västlig hakmossa",
    "Vågbandad barkbock"
) -join $crlf
$ws.Range("R4").Value = $artnamn4
